$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"

$ws.Range("C2").Select()
